$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column N (2020 data) ----------------------------------------

# N3: blank cell, bottom border only (same formatting family as M3)
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)   # xlPasteFormats

# N4: year header 2020 (same formatting as M4)
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N4").Value() = 2020

# N5: data value for 2020 (1.6) - uses a distinct style (Times New Roman 9,
# top+bottom medium border, vertically centered, general number format)
$ws.Range("A5").Copy()
$ws.Range("N5").PasteSpecial(-4122)   # xlPasteFormats
$n5 = $ws.Range("N5")
$n5.HorizontalAlignment = 1           # xlHAlignGeneral
$n5.WrapText = $false
$n5.Value() = 1.6

# --- Correct the 2019 figure -----------------------------------------
$ws.Range("M5").Value() = 1.6

# --- Restore selection as saved in the source file --------------------
[void]$ws.Range("P6").Select()

$excel.CutCopyMode = $false
